# Applies the commit:
#   - Refresh the cached "datetimeFigureOut" footer field from 3/2/22 to
#     3/8/22 everywhere it is cached (slide master, every slide layout,
#     and the notes master).
#   - Fix the protospacer/target sequence textbox on slide 1: the DNA
#     run "CCGGCTTGCAAACTCTCGCTCTA" becomes the RNA-style
#     "CCGGCUUGCAAACUCUCGCUCUA" (T -> U), leaving the trailing "N" run
#     (different colour/formatting) untouched.
#
# NOTE: this PS runtime does not reliably bind *named* function
# parameters (e.g. "-shapes $x"), so helper functions below are called
# positionally.

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.HasTextFrame -and ($sh.TextFrame.TextRange.Text -ne $newText)) {
                $sh.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

# --- 1. Slide master ------------------------------------------------------
$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes "3/8/22"

# --- 2. Every slide layout off the master ---------------------------------
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes "3/8/22"
}

# --- 3. Notes master --------------------------------------------------------
$notesMaster = $p.NotesMaster
Set-DatePlaceholderText $notesMaster.Shapes "3/8/22"

# --- 4. Protospacer/target sequence textbox on slide 1 ----------------------
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq "CCGGCTTGCAAACTCTCGCTCTAN") {
            $origHeight = $sh.Height
            $origWidth = $sh.Width
            # Only touch the first run (23 chars); the trailing "N" run has
            # its own distinct formatting and must stay untouched.
            $run = $tr.Characters(1, 23)
            $run.Text = "CCGGCUUGCAAACUCUCGCUCUA"
            # Editing the text range nudges the autosized box height by a
            # hair; put it back exactly as it was.
            $sh.Height = $origHeight
            $sh.Width = $origWidth
        }
    }
}
